$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '331.54'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '0.97%'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '41.47'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '3.99%'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.740'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '0.17%'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.08114'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '0.38%'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '8.671'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '-0.16%'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '4.499'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '-1.49%'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.972'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '0.26%'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.998'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '0.67%'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9272'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '-1.78%'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.1266'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '-1.16%'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.1958'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '-1.58%'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '8.783'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '15.52%'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.09191'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '-0.36%'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.03736'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '9.18%'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.1049'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '9.10%'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.001297'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '-1.35%'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.006337'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '4.00%'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.368'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '-0.11%'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '-0.98%'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.1368'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '-3.02%'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.2601'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '7.14%'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.04410'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '-0.75%'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.001252'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '-0.06%'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.004439'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '2.71%'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0001236'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '3.74%'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.02779'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '10.06%'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.05567'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '6.71%'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.007528'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '2.30%'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.009811'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '8.70%'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1423'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '-0.48%'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.002101'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '-4.17%'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.01185'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '18.31%'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00006759'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '1.12%'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.00000000747'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '-0.43%'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.003063'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '6.52%'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.002271'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '26.07%'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.00002092'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '-0.43%'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0001993'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '-0.43%'
